# The edit reorders data rows 2-20 on sheet "Artfynd":
#   - old row 2  -> new row 8
#   - old row 3  -> new row 20
#   - old rows 4..9   shift up to rows 2..7
#   - old rows 10..20 shift up to rows 9..19
#   - row 1 (header) and row 21 are unchanged
#
# This is implemented as a full-row move using two scratch rows (far below
# the used range) to stage rows 2 and 3 while the rest of the block shifts
# up, then those staged rows are dropped back in at their new positions.
#
# Each destination row is explicitly cleared before pasting into it because
# Copy()/paste here does not blank out destination cells that correspond to
# "empty" source cells (cells present in the XML with no value) - without
# the Clear() the old content would bleed through into columns that should
# end up empty (e.g. K, J, M, AC on various rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RowRange($r) {
    return $ws.Range("A" + $r + ":AY" + $r)
}

function MoveRow($srcRow, $dstRow) {
    (RowRange $dstRow).Clear()
    (RowRange $srcRow).Copy((RowRange $dstRow))
}

$scratch1 = 1000
$scratch2 = 1001

# Stage old row 2 and old row 3 out of the way.
(RowRange $scratch1).Clear()
(RowRange 2).Copy((RowRange $scratch1))
(RowRange $scratch2).Clear()
(RowRange 3).Copy((RowRange $scratch2))

# Shift old rows 4..9 up into 2..7.
MoveRow 4 2
MoveRow 5 3
MoveRow 6 4
MoveRow 7 5
MoveRow 8 6
MoveRow 9 7

# Drop staged old row 2 into its new spot at row 8.
MoveRow $scratch1 8

# Shift old rows 10..20 up into 9..19.
MoveRow 10 9
MoveRow 11 10
MoveRow 12 11
MoveRow 13 12
MoveRow 14 13
MoveRow 15 14
MoveRow 16 15
MoveRow 17 16
MoveRow 18 17
MoveRow 19 18
MoveRow 20 19

# Drop staged old row 3 into its new spot at row 20.
MoveRow $scratch2 20

# Clean up the scratch rows so nothing leaks outside the used range.
(RowRange $scratch1).Clear()
(RowRange $scratch2).Clear()
